$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.664.63"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.681.59"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.39"
$ws.Range("E4").Value = "  +26.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.999"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "228.70"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "652.20"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.441"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.14"
$ws.Range("E9").Value = "  +7.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.676.84"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.64"
$ws.Range("E12").Value = "  +7.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.208"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000300"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.56"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.365.51"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.403.52"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.88"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.678.67"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.83"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.83"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.539"
$ws.Range("E22").Value = "  +6.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "529.98"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.32"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.241"
$ws.Range("E25").Value = "  +41.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "119.57"
$ws.Range("E26").Value = "  +18.54%  "
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.81"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.878.33"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.84"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.29"
$ws.Range("E31").Value = "  +10.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.99"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.10"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.615"
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "611.97"
$ws.Range("E39").Value = "  -6.73%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.38"
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.04"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.163"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  +12.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.38"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.483"
$ws.Range("E46").Value = "  +11.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.955"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").Value = "  +5.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.54"
$ws.Range("E51").Value = "  -0.10%  "
